# Excel IC file edit - added TARV
# For each trace worksheet in the workbook, the Trace1 column (column B,
# rows 2-37) is updated to the TARV-adjusted values that were already
# staged in the Trace4 column (column E, rows 2-37) for simulations
# starting in August/September 2020.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    for ($r = 2; $r -le 37; $r++) {
        $newVal = $ws.Cells.Item($r, 5).Value2   # column E ("Trace4")
        $ws.Cells.Item($r, 2).Value2 = $newVal   # column B ("Trace1")
    }
}
